$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")
Write-Host $ws.Name
